# Apply the "proporcionalidades" scenario-coherence fix:
# - Remove the trailing scenario block (columns Y:AA)
# - Relabel the scenario headers in row 1/2 that were incoherent
# - Update the U:X data block (rows 3-6) to the corrected proportionalities

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the now-unneeded trailing columns Y:AA (whole columns)
$ws.Range("Y1:AA6").EntireColumn.Delete()

# 2. Fix row 1 scenario marker that used to live in column Y, now in W
$ws.Range("W1").Value = "6-2015"

# 3. Fix row 2 sub-header labels
$ws.Range("B2").Value = "3-2015"
$ws.Range("C2").Value = "4-2015"
$ws.Range("D2").Value = "5-2030-P"
$ws.Range("E2").Value = "3-2015"
$ws.Range("F2").Value = "5-2050-O"
$ws.Range("G2").Value = "4-2015"
$ws.Range("H2").Value = "5-2050-P"
$ws.Range("I2").Value = "3-2015"
$ws.Range("J2").Value = "5-2050-O"
$ws.Range("K2").Value = "4-2015"
$ws.Range("L2").Value = "5-2050-P"
$ws.Range("M2").Value = "3-2015"
$ws.Range("N2").Value = "5-2050-O"
$ws.Range("O2").Value = "4-2015"
$ws.Range("P2").Value = "5-2050-P"
$ws.Range("Q2").Value = "3-2015"
$ws.Range("R2").Value = "5-2030-O"
$ws.Range("S2").Value = "4-2015"
$ws.Range("T2").Value = "5-2015"
$ws.Range("U2").Value = "6-2015"
$ws.Range("V2").Value = "7-2015"
$ws.Range("W2").Value = "8-2015"
$ws.Range("X2").Value = "9-2015"

# 4. Update the corrected proportionality values for rows 3-6, columns U:X
$ws.Range("U3").Value = 0.738
$ws.Range("V3").Value = 0.262
$ws.Range("W3").Value = 0.442
$ws.Range("X3").Value = 0.558

$ws.Range("U4").Value = 0.738
$ws.Range("V4").Value = 0.262
$ws.Range("W4").Value = 0.442
$ws.Range("X4").Value = 0.558

$ws.Range("U5").Value = 0.738
$ws.Range("V5").Value = 0.262
$ws.Range("W5").Value = 0.442
$ws.Range("X5").Value = 0.558

$ws.Range("U6").Value = 0.738
$ws.Range("V6").Value = 0.262
$ws.Range("W6").Value = 0.442
$ws.Range("X6").Value = 0.558

# 5. Restore the selection/view state
$ws.Range("X2").Select() | Out-Null
